$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update capacity column (E) values for rows 2 and 3
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3

# Add new "shift" column in H
$ws.Range("H1").Value = "shift"
$ws.Range("H2").Value = 40
$ws.Range("H3").Value = 50

# Rename "is negative" header (currently in G1) to "isnegative"
$ws.Range("G1").Value = "isnegative"

# Update the active selection to match the target state
$ws.Range("G2").Select()
